# Camposanto.xlsx - aggiornamento fino a 27/05
# Append 14 new daily rows (r=256..269, dates 44330..44343) to Sheet1,
# mirroring the style/format of the existing data rows (col A uses the
# "date" cell style, columns B:D are plain numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$newRows = @(
    @(44330, 0, 0, 0),
    @(44331, 0, 0, 0),
    @(44332, 0, 0, 0),
    @(44333, 0, 0, 0),
    @(44334, 0, 0, 0),
    @(44335, 0, 0, 0),
    @(44336, 0, 0, 0),
    @(44337, 0, 0, 0),
    @(44338, 0, 0, 0),
    @(44339, 0, 0, 0),
    @(44340, 0, 0, 0),
    @(44341, 0, 0, 0),
    @(44342, 0, 0, 0),
    @(44343, 1, 1, 31.25976867771178)
)

$startRow = 256

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Copy the format of the last existing data row (255) down onto the
    # new row before writing values, so col A keeps the bordered/bold/
    # centered date style and B:D keep the plain default style.
    $ws.Range("A255:D255").Copy()
    $ws.Range("A$r:D$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}

$excel.CutCopyMode = 0
